$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 13213
$ws.Range("E2").Value = 1271
$ws.Range("F2").Value = 1271
$ws.Range("G2").Value = 1178
$ws.Range("H2").Value = 952
$ws.Range("I2").Value = 915
$ws.Range("J2").Value = 37
$ws.Range("K2").Value = 21224
$ws.Range("L2").Value = 6820
$ws.Range("M2").Value = 14403
$ws.Range("N2").Value = 13862
$ws.Range("O2").Value = 541
$ws.Range("P2").Value = 377
$ws.Range("Q2").Value = 1773
$ws.Range("R2").Value = 116
$ws.Range("S2").Value = -1978
$ws.Range("T2").Value = 565
$ws.Range("U2").Value = 1208
$ws.Range("V2").Value = 3138
$ws.Range("W2").Value = 9.619999999999999
$ws.Range("X2").Value = 7.2
$ws.Range("Y2").Value = 6.86
$ws.Range("Z2").Value = 4.46
$ws.Range("AA2").Value = 47.35
$ws.Range("AB2").Value = 3526.12
$ws.Range("AC2").Value = 12124
$ws.Range("AD2").Value = 12.63
$ws.Range("AE2").Value = 199834
$ws.Range("AF2").Value = 0.77
$ws.Range("AG2").Value = 1500
$ws.Range("AH2").Value = 0.98
$ws.Range("AI2").Value = 11.37
$ws.Range("AJ2").Value = 7545313

# Row 3
$ws.Range("D3").Value = 13773
$ws.Range("E3").Value = 1206
$ws.Range("F3").Value = 1206
$ws.Range("G3").Value = 1143
$ws.Range("H3").Value = 649
$ws.Range("I3").Value = 576
$ws.Range("J3").Value = 73
$ws.Range("K3").Value = 21507
$ws.Range("L3").Value = 6820
$ws.Range("M3").Value = 14687
$ws.Range("N3").Value = 14328
$ws.Range("O3").Value = 359
$ws.Range("P3").Value = 377
$ws.Range("Q3").Value = 1288
$ws.Range("R3").Value = 753
$ws.Range("S3").Value = -724
$ws.Range("T3").Value = 260
$ws.Range("U3").Value = 1028
$ws.Range("V3").Value = 2447
$ws.Range("W3").Value = 8.76
$ws.Range("X3").Value = 4.71
$ws.Range("Y3").Value = 4.09
$ws.Range("Z3").Value = 3.04
$ws.Range("AA3").Value = 46.43
$ws.Range("AB3").Value = 3641.46
$ws.Range("AC3").Value = 7633
$ws.Range("AD3").Value = 14.54
$ws.Range("AE3").Value = 206546
$ws.Range("AF3").Value = 0.54
$ws.Range("AG3").Value = 1500
$ws.Range("AH3").Value = 1.35
$ws.Range("AI3").Value = 18.07
$ws.Range("AJ3").Value = 7545313

# Row 4
$ws.Range("D4").Value = 14412
$ws.Range("E4").Value = 1016
$ws.Range("F4").Value = 1016
$ws.Range("G4").Value = 1301
$ws.Range("H4").Value = 838
$ws.Range("I4").Value = 808
$ws.Range("J4").Value = 29
$ws.Range("K4").Value = 21585
$ws.Range("L4").Value = 6304
$ws.Range("M4").Value = 15280
$ws.Range("N4").Value = 14845
$ws.Range("O4").Value = 435
$ws.Range("P4").Value = 377
$ws.Range("Q4").Value = 405
$ws.Range("R4").Value = 446
$ws.Range("S4").Value = -573
$ws.Range("T4").Value = 954
$ws.Range("U4").Value = -549
$ws.Range("V4").Value = 2039
$ws.Range("W4").Value = 7.05
$ws.Range("X4").Value = 5.81
$ws.Range("Y4").Value = 5.54
$ws.Range("Z4").Value = 3.89
$ws.Range("AA4").Value = 41.26
$ws.Range("AB4").Value = 3833.51
$ws.Range("AC4").Value = 10713
$ws.Range("AD4").Value = 7.46
$ws.Range("AE4").Value = 214002
$ws.Range("AF4").Value = 0.37
$ws.Range("AG4").Value = 1600
$ws.Range("AH4").Value = 2
$ws.Range("AI4").Value = 13.73
$ws.Range("AJ4").Value = 7545313

# Row 5
$ws.Range("D5").Value = 11908
$ws.Range("E5").Value = 1126
$ws.Range("F5").Value = 1126
$ws.Range("G5").Value = 864
$ws.Range("H5").Value = 656
$ws.Range("I5").Value = 571
$ws.Range("J5").Value = 85
$ws.Range("K5").Value = 33363
$ws.Range("L5").Value = 16969
$ws.Range("M5").Value = 16394
$ws.Range("N5").Value = 15337
$ws.Range("O5").Value = 1057
$ws.Range("P5").Value = 377
$ws.Range("Q5").Value = 704
$ws.Range("R5").Value = -8791
$ws.Range("S5").Value = 7195
$ws.Range("T5").Value = 585
$ws.Range("U5").Value = 119
$ws.Range("V5").Value = 10293
$ws.Range("W5").Value = 9.460000000000001
$ws.Range("X5").Value = 5.51
$ws.Range("Y5").Value = 3.79
$ws.Range("Z5").Value = 2.39
$ws.Range("AA5").Value = 103.5
$ws.Range("AB5").Value = 3952.12
$ws.Range("AC5").Value = 7574
$ws.Range("AD5").Value = 19.87
$ws.Range("AE5").Value = 221095
$ws.Range("AF5").Value = 0.68
$ws.Range("AG5").Value = 1800
$ws.Range("AH5").Value = 1.2
$ws.Range("AI5").Value = 21.85
$ws.Range("AJ5").Value = 7545313

# Row 6
$ws.Range("D6").Value = 11867
$ws.Range("E6").Value = 411
$ws.Range("F6").Value = 411
$ws.Range("G6").Value = 1378
$ws.Range("H6").Value = 5204
$ws.Range("I6").Value = 5117
$ws.Range("K6").Value = 32542
$ws.Range("L6").Value = 15294
$ws.Range("M6").Value = 17248
$ws.Range("N6").Value = 10002
$ws.Range("P6").Value = 314
$ws.Range("Q6").Value = 1267
$ws.Range("R6").Value = 2954
$ws.Range("S6").Value = -3343
$ws.Range("T6").Value = 691
$ws.Range("U6").Value = 576
$ws.Range("V6").Value = 8400
$ws.Range("W6").Value = 3.46
$ws.Range("X6").Value = 43.86
$ws.Range("Y6").Value = 40.39
$ws.Range("Z6").Value = 15.79
$ws.Range("AA6").Value = 88.67
$ws.Range("AB6").Value = 6586.75
$ws.Range("AC6").Value = 86875
$ws.Range("AD6").Value = 0.6
$ws.Range("AE6").Value = 166854
$ws.Range("AF6").Value = 0.31
$ws.Range("AG6").Value = 2200
$ws.Range("AH6").Value = 4.22
$ws.Range("AI6").Value = 2.58
$ws.Range("AJ6").Value = 5995844

# Row 7: clear removed cells
$ws.Range("D7:AJ7").ClearContents()

# Row 8: clear removed cells
$ws.Range("D8:AJ8").ClearContents()

# Row 9: clear removed cells
$ws.Range("D9:AJ9").ClearContents()
